$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell B11 currently holds the text "R40". It needs to become the text "1"
# (still stored as a shared string, not a number) while keeping its
# existing cell style/format untouched.
#
# A plain Value/Formula assignment of a numeric-looking string like "1"
# gets auto-coerced to a number by Excel, and forcing text via
# NumberFormat = "@" stamps a brand new cell style onto the cell. Instead,
# write it as a TEXT() formula (guaranteed string result) and then use
# Copy + PasteSpecial(values) back onto the very same cell: this freezes
# the formula result into a literal string value without touching the
# cell's number format/style and without disturbing any other cell.
$cell = $ws.Range("B11")
$cell.Formula = '=TEXT(1,"0")'
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163, 0)  # xlPasteValues, xlPasteSpecialOperationNone

$wb.Save()
